$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Mango" at Feria Lagunitas de
# Puerto Montt. In the source sheet this shows up as a brand-new row
# inserted at row 109 (pushing the existing rows 109-236 down to 110-237).
$ws.Rows("109").Insert()

# Populate the newly inserted row 109 with the new record's data.
$ws.Range("A109").Value = 4
$ws.Range("B109").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C109").Value = "Los Lagos"
$ws.Range("D109").Value = 44799
$ws.Range("E109").Value = 10
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100108
$ws.Range("H109").Value = "Tropicales y subtropicales"
$ws.Range("I109").Value = 100108002
$ws.Range("J109").Value = "Mango"
$ws.Range("K109").Value = "Sin especificar"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 200
$ws.Range("N109").Value = 13000
$ws.Range("O109").Value = 14000
$ws.Range("P109").Value = 13500
$ws.Range("Q109").Value = "`$/bandeja 4 kilos"
$ws.Range("R109").Value = "Brasil"
$ws.Range("S109").Value = 3375
$ws.Range("T109").Value = 4
